$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.706.06"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.296.39"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "96.91"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").Value = "268.58"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").Value = "45.37"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "15.71"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").Value = "2.641.40"
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "2.294.62"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "43.712.06"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "72.11"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "2.51"
$ws.Range("E22").Value = "  +10.25%  "
$ws.Range("D23").Value = "233.01"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("E25").Value = "  +5.93%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "11.28"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "38.76"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "175.10"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "21.86"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("D33").Value = "0.0905"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  +3.42%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "12.20"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.35"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").Value = "64.41"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("E46").Value = "  -5.04%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "97.27"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.53"
$ws.Range("E50").Value = "  +13.23%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.439"
$ws.Range("E51").Value = "  +4.65%  "
